# Test Data Added for Slovakia market
#
# Adds a new "Slovakia" worksheet (cloned from "Portugal", the previous
# last tab) right after "Portugal", fills in the market-specific cells,
# drops the "MZX Communicator" printer row that Slovakia doesn't use, and
# makes "Slovakia" the active/selected tab (Portugal goes back to a plain
# "select all" state, matching the rest of the non-active sheets).

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Capture the "whole sheet selected" state on Portugal while it is still
# the active sheet, before it stops being the active tab.
$portugal.Range("A1:XFD1048576").Select()

# Clone Portugal (this also places the copy immediately after it and
# makes the new sheet the active one, carrying over the selection above).
$portugal.Copy([System.Reflection.Missing]::Value, $portugal)

$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# Slovakia's printer list doesn't include "MZX Communicator" -> remove it.
$slovakia.Rows.Item(12).Delete()

# Market-specific cells.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3177/T3176/T3179"

# Rows 3 & 4 no longer need the taller wrapped-text row height used on the
# Portugal sheet; reset them back to the sheet's default height.
$slovakia.Range("A3:D4").EntireRow.AutoFit()

# Final selection on the new sheet.
$slovakia.Range("A8:A13").Select()
